$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values (FlashScore weekly games export) for rows 2-4

# Row 2 (Ind. Medellin - Alianza)
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 5.25
$ws.Range("J2").Value = 2.3
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("Z2").Value = 12
$ws.Range("AC2").Value = 8
$ws.Range("AO2").Value = 9
$ws.Range("AT2").Value = 2.5

# Row 3 (Santa Fe - Chico)
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62

# Row 4 (Correcaminos - Atl. Morelia)
$ws.Range("G4").Value = 2.6
$ws.Range("H4").Value = 3.05
$ws.Range("I4").Value = 2.65
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3.2
$ws.Range("M4").Value = 8.300000000000001
$ws.Range("N4").Value = 1.05
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 1.91
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 1.42
$ws.Range("T4").Value = 2.47
$ws.Range("U4").Value = 1.65
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 9.75
$ws.Range("Z4").Value = 30
$ws.Range("AA4").Value = 22
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 9.25
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 12.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 400
$ws.Range("AI4").Value = 14.5
$ws.Range("AJ4").Value = 9.5
$ws.Range("AK4").Value = 32
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 27
$ws.Range("AN4").Value = 4.55
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 23
$ws.Range("AQ4").Value = 70
$ws.Range("AS4").Value = 300
$ws.Range("AT4").Value = 2.42
$ws.Range("AU4").Value = 6.7
$ws.Range("AV4").Value = 60
$ws.Range("AW4").Value = 4.55
$ws.Range("AX4").Value = 14
$ws.Range("AY4").Value = 21
$ws.Range("AZ4").Value = 60
$ws.Range("BA4").Value = 90
$ws.Range("BB4").Value = 200
